# "Delay in Message recieved"
# The workbook originally had three sheets: "Binance", "Binance Timing" and
# "Binance Occ.". Only "Binance Timing" survives, and it is populated with a
# small table of message-delay timing data (a pair/trio of timestamps plus
# the elapsed-time deltas between them).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Drop the two sheets that are no longer needed -------------------------
$wb.Worksheets.Item("Binance").Delete()
$wb.Worksheets.Item("Binance Occ.").Delete()

$ws = $wb.Worksheets.Item("Binance Timing")
$ws.Activate()

# --- Column widths -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 17.5
$ws.Columns.Item(2).ColumnWidth = 19.5
$ws.Columns.Item(3).ColumnWidth = 18.166666666666668
$ws.Columns.Item(4).ColumnWidth = 21.666666666666668

# --- Header-ish row: the currency pairs / trade leg labels ------------------
$ws.Range("A2").Value = "buy"
$ws.Range("B2").Value = "AUDUSDT"
$ws.Range("C2").Value = "LUNAAUD"
$ws.Range("D2").Value = "LUNAUSDT"

# --- Computed ratio ----------------------------------------------------
$ws.Range("A3").Value = 1.91838985784365

# --- Timestamps (yyyy-mm-dd h:mm:ss) and elapsed deltas ([hh]:mm:ss) --------
$ws.Range("A4").Value = 44693.6622156713
$ws.Range("A4").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("A5").Value = 44693.66239370311
$ws.Range("A5").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("A6").Value = 0.0001780318171296296
$ws.Range("A6").NumberFormat = "[hh]:mm:ss"

$ws.Range("A7").Value = 44693.66239580378
$ws.Range("A7").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("A8").Value = 0.000002100671296296296
$ws.Range("A8").NumberFormat = "[hh]:mm:ss"

# --- Cosmetics: selection + print setup matching the authored file ---------
$ws.Range("C36").Select()
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9

Write-Output "Binance Timing populated; Binance / Binance Occ. removed"
